$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Cyprinidae"
$ws.Range("B2").Value = "Nemacheilidae"
$ws.Range("C2").Value = 82
$ws.Range("D2").Value = 999
$ws.Range("E2").Value = 5.614677304204756
$ws.Range("F2").Value = 0.001
$ws.Range("G2").Value = 0.005
$ws.Range("H2").Value = "skin"

$ws.Range("A3").Value = "Nemacheilidae"
$ws.Range("B3").Value = "Tilapiinae"
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 999
$ws.Range("E3").Value = 5.25720210706947
$ws.Range("F3").Value = 0.001
$ws.Range("G3").Value = 0.005
$ws.Range("H3").Value = "skin"

$ws.Range("A4").Value = "Haplochrominae"
$ws.Range("B4").Value = "Nemacheilidae"
$ws.Range("C4").Value = 21
$ws.Range("D4").Value = 999
$ws.Range("E4").Value = 3.946261484657737
$ws.Range("F4").Value = 0.002
$ws.Range("G4").Value = 0.006666666666666667
$ws.Range("H4").Value = "skin"

$ws.Range("A5").Value = "Cyprinidae"
$ws.Range("B5").Value = "Tilapiinae"
$ws.Range("C5").Value = 76
$ws.Range("D5").Value = 999
$ws.Range("E5").Value = 2.199024440688461
$ws.Range("F5").Value = 0.015
$ws.Range("G5").Value = 0.0375
$ws.Range("H5").Value = "skin"

$ws.Range("A6").Value = "Cyprinidae"
$ws.Range("B6").Value = "Haplochrominae"
$ws.Range("C6").Value = 67
$ws.Range("D6").Value = 999
$ws.Range("E6").Value = 1.719444743503678
$ws.Range("F6").Value = 0.05
$ws.Range("G6").Value = 0.08333333333333334
$ws.Range("H6").Value = "skin"

$ws.Range("A7").Value = "Mugilidae"
$ws.Range("B7").Value = "Nemacheilidae"
$ws.Range("C7").Value = 19
$ws.Range("D7").Value = 999
$ws.Range("E7").Value = 2.565531787199989
$ws.Range("F7").Value = 0.044
$ws.Range("G7").Value = 0.08333333333333334
$ws.Range("H7").Value = "skin"

$ws.Range("A8").Value = "Cyprinidae"
$ws.Range("B8").Value = "Mugilidae"
$ws.Range("C8").Value = 65
$ws.Range("D8").Value = 999
$ws.Range("E8").Value = 1.413268315735274
$ws.Range("F8").Value = 0.096
$ws.Range("G8").Value = 0.1371428571428571
$ws.Range("H8").Value = "skin"

$ws.Range("A9").Value = "Haplochrominae"
$ws.Range("B9").Value = "Tilapiinae"
$ws.Range("C9").Value = 15
$ws.Range("D9").Value = 999
$ws.Range("E9").Value = 0.9193645798307706
$ws.Range("F9").Value = 0.526
$ws.Range("G9").Value = 0.5844444444444444
$ws.Range("H9").Value = "skin"

$ws.Range("A10").Value = "Mugilidae"
$ws.Range("B10").Value = "Tilapiinae"
$ws.Range("C10").Value = 13
$ws.Range("D10").Value = 999
$ws.Range("E10").Value = 1.01936922172975
$ws.Range("F10").Value = 0.47
$ws.Range("G10").Value = 0.5844444444444444
$ws.Range("H10").Value = "skin"

$ws.Range("A11").Value = "Haplochrominae"
$ws.Range("B11").Value = "Mugilidae"
$ws.Range("C11").Value = 4
$ws.Range("D11").Value = 999
$ws.Range("E11").Value = 0.4578020413090889
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = "skin"

$ws.Range("A12").Value = "Cyprinidae"
$ws.Range("B12").Value = "Haplochrominae"
$ws.Range("C12").Value = 72
$ws.Range("D12").Value = 999
$ws.Range("E12").Value = 4.221830888078177
$ws.Range("F12").Value = 0.001
$ws.Range("G12").Value = 0.003
$ws.Range("H12").Value = "swab"

$ws.Range("A13").Value = "Cyprinidae"
$ws.Range("B13").Value = "Nemacheilidae"
$ws.Range("C13").Value = 85
$ws.Range("D13").Value = 999
$ws.Range("E13").Value = 3.247701309998819
$ws.Range("F13").Value = 0.001
$ws.Range("G13").Value = 0.003
$ws.Range("H13").Value = "swab"

$ws.Range("A14").Value = "Cyprinidae"
$ws.Range("B14").Value = "Tilapiinae"
$ws.Range("C14").Value = 86
$ws.Range("D14").Value = 999
$ws.Range("E14").Value = 5.551815377807659
$ws.Range("F14").Value = 0.001
$ws.Range("G14").Value = 0.003
$ws.Range("H14").Value = "swab"

$ws.Range("A15").Value = "Haplochrominae"
$ws.Range("B15").Value = "Nemacheilidae"
$ws.Range("C15").Value = 29
$ws.Range("D15").Value = 999
$ws.Range("E15").Value = 4.198163177854015
$ws.Range("F15").Value = 0.001
$ws.Range("G15").Value = 0.003
$ws.Range("H15").Value = "swab"

$ws.Range("A16").Value = "Nemacheilidae"
$ws.Range("B16").Value = "Tilapiinae"
$ws.Range("C16").Value = 43
$ws.Range("D16").Value = 999
$ws.Range("E16").Value = 4.572618327401713
$ws.Range("F16").Value = 0.001
$ws.Range("G16").Value = 0.003
$ws.Range("H16").Value = "swab"

$ws.Range("A17").Value = "Cyprinidae"
$ws.Range("B17").Value = "Mugilidae"
$ws.Range("C17").Value = 67
$ws.Range("D17").Value = 999
$ws.Range("E17").Value = 2.678114025269259
$ws.Range("F17").Value = 0.002
$ws.Range("G17").Value = 0.004285714285714286
$ws.Range("H17").Value = "swab"

$ws.Range("A18").Value = "Cyprinidae"
$ws.Range("B18").Value = "Poeciliidae"
$ws.Range("C18").Value = 66
$ws.Range("D18").Value = 999
$ws.Range("E18").Value = 2.159926576250572
$ws.Range("F18").Value = 0.002
$ws.Range("G18").Value = 0.004285714285714286
$ws.Range("H18").Value = "swab"

$ws.Range("A19").Value = "Mugilidae"
$ws.Range("B19").Value = "Nemacheilidae"
$ws.Range("C19").Value = 24
$ws.Range("D19").Value = 999
$ws.Range("E19").Value = 2.848382228920664
$ws.Range("F19").Value = 0.004
$ws.Range("G19").Value = 0.007500000000000001
$ws.Range("H19").Value = "swab"

$ws.Range("A20").Value = "Nemacheilidae"
$ws.Range("B20").Value = "Poeciliidae"
$ws.Range("C20").Value = 23
$ws.Range("D20").Value = 999
$ws.Range("E20").Value = 2.34527189141601
$ws.Range("F20").Value = 0.019
$ws.Range("G20").Value = 0.03166666666666667
$ws.Range("H20").Value = "swab"

$ws.Range("A21").Value = "Mugilidae"
$ws.Range("B21").Value = "Poeciliidae"
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 999
$ws.Range("E21").Value = 1.594030007225507
$ws.Range("F21").Value = 0.11
$ws.Range("G21").Value = 0.165
$ws.Range("H21").Value = "swab"

$ws.Range("A22").Value = "Haplochrominae"
$ws.Range("B22").Value = "Poeciliidae"
$ws.Range("C22").Value = 10
$ws.Range("D22").Value = 999
$ws.Range("E22").Value = 1.172278478121282
$ws.Range("F22").Value = 0.135
$ws.Range("G22").Value = 0.1840909090909091
$ws.Range("H22").Value = "swab"

$ws.Range("A23").Value = "Poeciliidae"
$ws.Range("B23").Value = "Tilapiinae"
$ws.Range("C23").Value = 24
$ws.Range("D23").Value = 999
$ws.Range("E23").Value = 1.132092876578544
$ws.Range("F23").Value = 0.264
$ws.Range("G23").Value = 0.33
$ws.Range("H23").Value = "swab"

$ws.Range("A24").Value = "Mugilidae"
$ws.Range("B24").Value = "Tilapiinae"
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 999
$ws.Range("E24").Value = 1.086555316335504
$ws.Range("F24").Value = 0.328
$ws.Range("G24").Value = 0.3784615384615385
$ws.Range("H24").Value = "swab"

$ws.Range("A25").Value = "Haplochrominae"
$ws.Range("B25").Value = "Mugilidae"
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 999
$ws.Range("E25").Value = 0.9997068112313988
$ws.Range("F25").Value = 0.385
$ws.Range("G25").Value = 0.4125
$ws.Range("H25").Value = "swab"

$ws.Range("A26").Value = "Haplochrominae"
$ws.Range("B26").Value = "Tilapiinae"
$ws.Range("C26").Value = 30
$ws.Range("D26").Value = 999
$ws.Range("E26").Value = 0.8089070422353125
$ws.Range("F26").Value = 0.684
$ws.Range("G26").Value = 0.684
$ws.Range("H26").Value = "swab"

